# Auto-generated edit script for 南宁-漫展信息.xlsx
# Applies: (1) sheet1 (展览) F-value updates for 3 events;
#          (2) sheet2 (演出) drop the 2024-08-10 row, shift remaining rows up, delete trailing row;
#          (3) sheet4 (全部类型) same drop/shift/delete pattern as sheet2.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (index 1): update 想去人数 (F) for the three listed events ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 1293
$ws1.Range("F3").Value = 2798
$ws1.Range("F4").Value = 250

# --- Sheet "演出" (index 2): the 2024-08-10 event is gone; remaining events shift up ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = '2024-08-14'
$ws2.Range("C2").Value = '南宁·新西兰治愈系民谣歌手LukeThompson2024中国巡演 KEEP ROLLING ON '
$ws2.Range("D2").Value = '中山路万象汇L2层37号 候朋现场HOPELIVE-中山路万象汇店'
$ws2.Range("E2").Value = '2024.08.14 20:00-08.14 21:30'
$ws2.Range("F2").Value = 4
$ws2.Range("G2").Value = 180
$ws2.Range("H2").Value = 'https://show.bilibili.com/platform/detail.html?id=88015'
$ws2.Range("I2").Value = '//i1.hdslb.com/bfs/openplatform/202406/76WI4tA01718179482365.jpeg'
$ws2.Range("B3").Value = '2024-10-04'
$ws2.Range("C3").Value = '南宁·《最后的莫西干人——亚历桑德罗&丛林回响乐队印第安音乐品鉴会》'
$ws2.Range("D3").Value = '福建园街道星光大道4号 南宁剧场'
$ws2.Range("E3").Value = '2024.10.04 20:00-10.04 21:30'
$ws2.Range("F3").Value = 6
$ws2.Range("G3").Value = 100
$ws2.Range("H3").Value = 'https://show.bilibili.com/platform/detail.html?id=89039'
$ws2.Range("I3").Value = '//i0.hdslb.com/bfs/openplatform/202407/dudapgjU1720595605665.jpeg'
$ws2.Rows.Item(4).Delete()

# --- Sheet "全部类型" (index 4): same 2024-08-10 removal, rows shift up ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = '2024-08-14'
$ws4.Range("C2").Value = '南宁·新西兰治愈系民谣歌手LukeThompson2024中国巡演 KEEP ROLLING ON '
$ws4.Range("D2").Value = '中山路万象汇L2层37号 候朋现场HOPELIVE-中山路万象汇店'
$ws4.Range("E2").Value = '2024.08.14 20:00-08.14 21:30'
$ws4.Range("F2").Value = 4
$ws4.Range("G2").Value = 180
$ws4.Range("H2").Value = 'https://show.bilibili.com/platform/detail.html?id=88015'
$ws4.Range("I2").Value = '//i1.hdslb.com/bfs/openplatform/202406/76WI4tA01718179482365.jpeg'
$ws4.Range("B3").Value = '2024-08-17'
$ws4.Range("C3").Value = '南宁·熊喵M动漫嘉年华【免费】'
$ws4.Range("D3").Value = '港航上尧码头(江北大道南100米) 水明漾艺术中心'
$ws4.Range("E3").Value = '2024.08.17 10:00-08.18 17:00'
$ws4.Range("F3").Value = 1293
$ws4.Range("G3").Value = 29.9
$ws4.Range("H3").Value = 'https://show.bilibili.com/platform/detail.html?id=89145'
$ws4.Range("I3").Value = '//i2.hdslb.com/bfs/openplatform/202407/ndmB7MOh1720344131003.jpeg'
$ws4.Range("B4").Value = '2024-08-24'
$ws4.Range("C4").Value = '南宁·第二届北极光动漫展'
$ws4.Range("D4").Value = '民族大道106号 南宁国际会展中心'
$ws4.Range("E4").Value = '2024.08.24 09:00-08.25 17:00'
$ws4.Range("F4").Value = 2798
$ws4.Range("G4").Value = 65
$ws4.Range("H4").Value = 'https://show.bilibili.com/platform/detail.html?id=88276'
$ws4.Range("I4").Value = '//i1.hdslb.com/bfs/openplatform/202406/mTEwC1GY1717576221099.jpeg'
$ws4.Range("B5").Value = '2024-10-04'
$ws4.Range("C5").Value = '南宁·《最后的莫西干人——亚历桑德罗&丛林回响乐队印第安音乐品鉴会》'
$ws4.Range("D5").Value = '福建园街道星光大道4号 南宁剧场'
$ws4.Range("E5").Value = '2024.10.04 20:00-10.04 21:30'
$ws4.Range("F5").Value = 6
$ws4.Range("G5").Value = 100
$ws4.Range("H5").Value = 'https://show.bilibili.com/platform/detail.html?id=89039'
$ws4.Range("I5").Value = '//i0.hdslb.com/bfs/openplatform/202407/dudapgjU1720595605665.jpeg'
$ws4.Range("B6").Value = '2024-11-02'
$ws4.Range("C6").Value = '南宁·万圣漫控嘉年华10'
$ws4.Range("D6").Value = '亭洪路45号 百益上河城'
$ws4.Range("E6").Value = '2024.11.02 11:00-11.03 22:00'
$ws4.Range("F6").Value = 250
$ws4.Range("G6").Value = 50
$ws4.Range("H6").Value = 'https://show.bilibili.com/platform/detail.html?id=87820'
$ws4.Range("I6").Value = '//i1.hdslb.com/bfs/openplatform/202406/abJD2cvV1718955681653.jpeg'
$ws4.Rows.Item(7).Delete()
